# Apply updated evaluation metrics across the three worksheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.7078651685393258
$wsSummary.Range("C2").Value = 0.937007874015748
$wsSummary.Range("D2").Value = 0.4456928838951311
$wsSummary.Range("E2").Value = 0.6040609137055838
$wsSummary.Range("F2").Value = 0.497907949790795
$wsSummary.Range("G2").Value = 0.4548662158188768
$wsSummary.Range("H2").Value = 0.7078651685393259
$wsSummary.Range("I2").Value = 238
$wsSummary.Range("J2").Value = 16
$wsSummary.Range("K2").Value = 518
$wsSummary.Range("L2").Value = 296

# --- Sheet 2: Classification Report ---
$wsReport = $wb.Worksheets.Item("Classification Report")
$wsReport.Range("B2").Value = 0.6363636363636364
$wsReport.Range("C2").Value = 0.9700374531835206
$wsReport.Range("D2").Value = 0.7685459940652819

$wsReport.Range("B3").Value = 0.937007874015748
$wsReport.Range("C3").Value = 0.4456928838951311
$wsReport.Range("D3").Value = 0.6040609137055838

$wsReport.Range("B4").Value = 0.7078651685393258
$wsReport.Range("C4").Value = 0.7078651685393258
$wsReport.Range("D4").Value = 0.7078651685393258
$wsReport.Range("E4").Value = 0.7078651685393258

$wsReport.Range("B5").Value = 0.7866857551896922
$wsReport.Range("C5").Value = 0.7078651685393258
$wsReport.Range("D5").Value = 0.6863034538854329

$wsReport.Range("B6").Value = 0.7866857551896922
$wsReport.Range("C6").Value = 0.7078651685393258
$wsReport.Range("D6").Value = 0.6863034538854329

# --- Sheet 3: Confusion Matrix ---
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")
$wsConfusion.Range("B2").Value = 518
$wsConfusion.Range("C2").Value = 16

$wsConfusion.Range("B3").Value = 296
$wsConfusion.Range("C3").Value = 238
